$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "ECs" sending-cluster rows (rows 2-4); remaining rows shift up.
$ws.Rows("2:4").Delete()

# Update remaining rows (now rows 2-7) with recalculated TPM-derived values.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Alcam"
$ws.Range("C2").Value = "Chl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5683613333333334
$ws.Range("H2").Value = 1.705084
$ws.Range("I2").Value = 0.4361027177196302
$ws.Range("J2").Value = 0.4361027177196302
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005250333333333333
$ws.Range("N2").Value = 0.015751
$ws.Range("O2").Value = 0.001037378148736751
$ws.Range("P2").Value = 0.001037378148736751
$ws.Range("Q2").Value = 0.002984086453777778
$ws.Range("R2").Value = 0.026856778084
$ws.Range("S2").Value = 0.000452403429967056
$ws.Range("T2").Value = 0.000452403429967056
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Alcam"
$ws.Range("C3").Value = "Chl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5683613333333334
$ws.Range("H3").Value = 1.705084
$ws.Range("I3").Value = 0.4361027177196302
$ws.Range("J3").Value = 0.4361027177196302
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.466778
$ws.Range("N3").Value = 1.400334
$ws.Range("O3").Value = 0.0922275342856409
$ws.Range("P3").Value = 0.0922275342856409
$ws.Range("Q3").Value = 0.2652985664506667
$ws.Range("R3").Value = 2.387687098056
$ws.Range("S3").Value = 0.04022067835054837
$ws.Range("T3").Value = 0.04022067835054837
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Alcam"
$ws.Range("C4").Value = "Chl1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5683613333333334
$ws.Range("H4").Value = 1.705084
$ws.Range("I4").Value = 0.4361027177196302
$ws.Range("J4").Value = 0.4361027177196302
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.589128333333333
$ws.Range("N4").Value = 13.767385
$ws.Range("O4").Value = 0.9067350875656224
$ws.Range("P4").Value = 0.9067350875656223
$ws.Range("Q4").Value = 2.608283098371111
$ws.Range("R4").Value = 23.47454788534
$ws.Range("S4").Value = 0.3954296359391148
$ws.Range("T4").Value = 0.3954296359391147
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Alcam"
$ws.Range("C5").Value = "Chl1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7349126666666667
$ws.Range("H5").Value = 2.204738
$ws.Range("I5").Value = 0.5638972822803697
$ws.Range("J5").Value = 0.5638972822803697
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.005250333333333333
$ws.Range("N5").Value = 0.015751
$ws.Range("O5").Value = 0.001037378148736751
$ws.Range("P5").Value = 0.001037378148736751
$ws.Range("Q5").Value = 0.003858536470888889
$ws.Range("R5").Value = 0.034726828238
$ws.Range("S5").Value = 0.0005849747187696952
$ws.Range("T5").Value = 0.0005849747187696952
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Alcam"
$ws.Range("C6").Value = "Chl1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7349126666666667
$ws.Range("H6").Value = 2.204738
$ws.Range("I6").Value = 0.5638972822803697
$ws.Range("J6").Value = 0.5638972822803697
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.466778
$ws.Range("N6").Value = 1.400334
$ws.Range("O6").Value = 0.0922275342856409
$ws.Range("P6").Value = 0.0922275342856409
$ws.Range("Q6").Value = 0.3430410647213333
$ws.Range("R6").Value = 3.087369582492
$ws.Range("S6").Value = 0.05200685593509252
$ws.Range("T6").Value = 0.05200685593509252
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Alcam"
$ws.Range("C7").Value = "Chl1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7349126666666667
$ws.Range("H7").Value = 2.204738
$ws.Range("I7").Value = 0.5638972822803697
$ws.Range("J7").Value = 0.5638972822803697
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.589128333333333
$ws.Range("N7").Value = 13.767385
$ws.Range("O7").Value = 0.9067350875656224
$ws.Range("P7").Value = 0.9067350875656223
$ws.Range("Q7").Value = 3.372608541125555
$ws.Range("R7").Value = 30.35347687013
$ws.Range("S7").Value = 0.5113054516265075
$ws.Range("T7").Value = 0.5113054516265074
